$d = $word.ActiveDocument

# The document's TOC-linked bookmarks (_Toc########) were regenerated by
# Word with new, higher numeric suffixes. Rename each one in place,
# preserving its position/range, by adding a bookmark with the new name
# over the old bookmark's range and then deleting the old one.
$bookmarkRenames = @(
    @{ Old = "_Toc24731381"; New = "_Toc25148497" },
    @{ Old = "_Toc24731382"; New = "_Toc25148498" },
    @{ Old = "_Toc24731383"; New = "_Toc25148499" },
    @{ Old = "_Toc24731384"; New = "_Toc25148500" },
    @{ Old = "_Toc24731385"; New = "_Toc25148501" },
    @{ Old = "_Toc24731386"; New = "_Toc25148502" },
    @{ Old = "_Toc24731387"; New = "_Toc25148503" },
    @{ Old = "_Toc24731388"; New = "_Toc25148504" },
    @{ Old = "_Toc24731389"; New = "_Toc25148505" },
    @{ Old = "_Toc24731390"; New = "_Toc25148506" },
    @{ Old = "_Toc24731391"; New = "_Toc25148507" },
    @{ Old = "_Toc24731392"; New = "_Toc25148508" },
    @{ Old = "_Toc24731393"; New = "_Toc25148509" },
    @{ Old = "_Toc24731394"; New = "_Toc25148510" },
    @{ Old = "_Toc24731395"; New = "_Toc25148511" },
    @{ Old = "_Toc24731396"; New = "_Toc25148512" },
    @{ Old = "_Toc24731397"; New = "_Toc25148513" },
    @{ Old = "_Toc24731398"; New = "_Toc25148514" },
    @{ Old = "_Toc24731399"; New = "_Toc25148515" }
)

foreach ($pair in $bookmarkRenames) {
    if ($d.Bookmarks.Exists($pair.Old)) {
        $bm = $d.Bookmarks.Item($pair.Old)
        $r = $bm.Range
        $d.Bookmarks.Add($pair.New, $r)
        $bm.Delete()
    }
}

# Remove the stray "_GoBack" bookmark (last-edit-position marker) that was
# left next to "Idera DB Optimizer".
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
